$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(28, 8).Value = 10994.5   # H28: 10008.546 -> 10994.5
$ws.Cells.Item(28, 9).Value = 14650.714   # I28: 12838 -> 14650.714
$ws.Cells.Item(28, 11).Value = 14650.714   # K28: 12838 -> 14650.714
$ws.Cells.Item(28, 13).Value = -14165.714   # M28: -12353 -> -14165.714
$ws.Cells.Item(87, 8).Value = 49999   # H87: 44999 -> 49999
$ws.Cells.Item(87, 10).Value = 49999   # J87: 44999 -> 49999
$ws.Cells.Item(87, 12).Value = 49999   # L87: 44999 -> 49999
$ws.Cells.Item(87, 14).Value = -52495   # N87: -47495 -> -52495
$ws.Cells.Item(90, 8).Value = 49999   # H90: 44999 -> 49999
$ws.Cells.Item(90, 10).Value = 49999   # J90: 44999 -> 49999
$ws.Cells.Item(90, 12).Value = 149997   # L90: 134997 -> 149997
$ws.Cells.Item(90, 14).Value = -162477   # N90: -147477 -> -162477
$ws.Cells.Item(95, 8).Value = 21641.334   # H95: 29974.334 -> 21641.334
$ws.Cells.Item(95, 10).Value = 21641.334   # J95: 29974.334 -> 21641.334
$ws.Cells.Item(95, 12).Value = 21641.334   # L95: 29974.334 -> 21641.334
$ws.Cells.Item(95, 14).Value = -27133.334   # N95: -35466.334 -> -27133.334
$ws.Cells.Item(100, 8).Value = 4099.6   # H100: 4549.8 -> 4099.6
$ws.Cells.Item(100, 9).Value = 4099.6   # I100: 4549.8 -> 4099.6
$ws.Cells.Item(100, 11).Value = 4099.6   # K100: 4549.8 -> 4099.6
$ws.Cells.Item(100, 13).Value = -3558.6   # M100: -4008.8 -> -3558.6
$ws.Cells.Item(106, 8).Value = 4778   # H106: 4914 -> 4778
$ws.Cells.Item(106, 9).Value = 4778   # I106: 4914 -> 4778
$ws.Cells.Item(106, 11).Value = 4778   # K106: 4914 -> 4778
$ws.Cells.Item(106, 13).Value = -4147   # M106: -4283 -> -4147
$ws.Cells.Item(135, 8).Value = 398.1   # H135: 434.55554 -> 398.1
$ws.Cells.Item(135, 9).Value = 199.5   # I135: 218 -> 199.5
$ws.Cells.Item(135, 11).Value = 1795.5   # K135: 1962 -> 1795.5
$ws.Cells.Item(135, 13).Value = 739.5   # M135: 573 -> 739.5
$ws.Cells.Item(138, 8).Value = 3769.5715   # H138: 3906.4614 -> 3769.5715
$ws.Cells.Item(138, 9).Value = 2129   # I138: 2156.8 -> 2129
$ws.Cells.Item(138, 11).Value = 6387   # K138: 6470.400000000001 -> 6387
$ws.Cells.Item(138, 13).Value = -1247   # M138: -1330.400000000001 -> -1247

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(19, 8).Value = 11736   # H19: 9027 -> 11736
$ws.Cells.Item(19, 10).Value = 17500   # J19: 11966.667 -> 17500
$ws.Cells.Item(19, 12).Value = 17500   # L19: 11966.667 -> 17500
$ws.Cells.Item(19, 14).Value = -17958   # N19: -12424.667 -> -17958
$ws.Cells.Item(61, 8).Value = 0   # H61: 3469.2856 -> 0
$ws.Cells.Item(61, 9).Value = 0   # I61: 3469.2856 -> 0
$ws.Cells.Item(61, 11).Value = 0   # K61: 3469.2856 -> 0
$ws.Cells.Item(61, 13).Value = $null  # M61: clear (was -3257.2856)
$ws.Cells.Item(74, 8).Value = 5462.364   # H74: 5202.2 -> 5462.364
$ws.Cells.Item(74, 9).Value = 4584   # I74: 4548.636 -> 4584
$ws.Cells.Item(74, 11).Value = 4584   # K74: 4548.636 -> 4584
$ws.Cells.Item(74, 13).Value = -3710   # M74: -3674.636 -> -3710
$ws.Cells.Item(77, 8).Value = 5462.364   # H77: 5202.2 -> 5462.364
$ws.Cells.Item(77, 9).Value = 4584   # I77: 4548.636 -> 4584
$ws.Cells.Item(77, 11).Value = 22920   # K77: 22743.18 -> 22920
$ws.Cells.Item(77, 13).Value = -18552   # M77: -18375.18 -> -18552
$ws.Cells.Item(97, 8).Value = 1043.2222   # H97: 1095.125 -> 1043.2222
$ws.Cells.Item(97, 9).Value = 999.7143   # I97: 1037.2858 -> 999.7143
$ws.Cells.Item(97, 10).Value = 1195.5   # J97: 1500 -> 1195.5
$ws.Cells.Item(97, 11).Value = 999.7143   # K97: 1037.2858 -> 999.7143
$ws.Cells.Item(97, 12).Value = 1195.5   # L97: 1500 -> 1195.5
$ws.Cells.Item(97, 13).Value = -503.7143   # M97: -541.2858000000001 -> -503.7143
$ws.Cells.Item(97, 14).Value = -2187.5   # N97: -2492 -> -2187.5
$ws.Cells.Item(102, 8).Value = 2100.6365   # H102: 2102.2727 -> 2100.6365
$ws.Cells.Item(102, 9).Value = 1311.2   # I102: 1313 -> 1311.2
$ws.Cells.Item(102, 11).Value = 1311.2   # K102: 1313 -> 1311.2
$ws.Cells.Item(102, 13).Value = 310.8   # M102: 309 -> 310.8
$ws.Cells.Item(110, 8).Value = 2334.4443   # H110: 1919.8462 -> 2334.4443
$ws.Cells.Item(110, 9).Value = 2288.5715   # I110: 1897.8 -> 2288.5715
$ws.Cells.Item(110, 10).Value = 2495   # J110: 1993.3334 -> 2495
$ws.Cells.Item(110, 11).Value = 2288.5715   # K110: 1897.8 -> 2288.5715
$ws.Cells.Item(110, 12).Value = 2495   # L110: 1993.3334 -> 2495
$ws.Cells.Item(110, 13).Value = -243.5715   # M110: 147.2 -> -243.5715
$ws.Cells.Item(110, 14).Value = -6585   # N110: -6083.3334 -> -6585
$ws.Cells.Item(132, 8).Value = 1633.5714   # H132: 1705.25 -> 1633.5714
$ws.Cells.Item(132, 9).Value = 1706.0526   # I132: 1789.7222 -> 1706.0526
$ws.Cells.Item(132, 11).Value = 5118.1578   # K132: 5369.1666 -> 5118.1578
$ws.Cells.Item(132, 13).Value = -2588.1578   # M132: -2839.1666 -> -2588.1578
$ws.Cells.Item(136, 8).Value = 0   # H136: 3469.2856 -> 0
$ws.Cells.Item(136, 9).Value = 0   # I136: 3469.2856 -> 0
$ws.Cells.Item(136, 11).Value = 0   # K136: 10407.8568 -> 0
$ws.Cells.Item(136, 13).Value = $null  # M136: clear (was -7857.856800000001)

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 4010.682   # H20: 4163.2383 -> 4010.682
$ws.Cells.Item(20, 9).Value = 3596.7334   # I20: 3796 -> 3596.7334
$ws.Cells.Item(20, 11).Value = 3596.7334   # K20: 3796 -> 3596.7334
$ws.Cells.Item(20, 13).Value = -3349.7334   # M20: -3549 -> -3349.7334
$ws.Cells.Item(105, 8).Value = 2335.3333   # H105: 1627.7273 -> 2335.3333
$ws.Cells.Item(105, 9).Value = 2003   # I105: 1393.5714 -> 2003
$ws.Cells.Item(105, 10).Value = 3000   # J105: 2037.5 -> 3000
$ws.Cells.Item(105, 11).Value = 2003   # K105: 1393.5714 -> 2003
$ws.Cells.Item(105, 12).Value = 3000   # L105: 2037.5 -> 3000
$ws.Cells.Item(105, 13).Value = -256   # M105: 353.4286 -> -256
$ws.Cells.Item(105, 14).Value = -6494   # N105: -5531.5 -> -6494
$ws.Cells.Item(107, 8).Value = 998.25   # H107: 999.6667 -> 998.25
$ws.Cells.Item(107, 9).Value = 998.25   # I107: 999.6667 -> 998.25
$ws.Cells.Item(107, 11).Value = 998.25   # K107: 999.6667 -> 998.25
$ws.Cells.Item(107, 13).Value = 921.75   # M107: 920.3333 -> 921.75

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 1242.75   # H16: 1191.4 -> 1242.75
$ws.Cells.Item(16, 9).Value = 1242.75   # I16: 1191.4 -> 1242.75
$ws.Cells.Item(16, 11).Value = 1242.75   # K16: 1191.4 -> 1242.75
$ws.Cells.Item(16, 13).Value = -955.75   # M16: -904.4000000000001 -> -955.75
$ws.Cells.Item(113, 8).Value = 1242.75   # H113: 1191.4 -> 1242.75
$ws.Cells.Item(113, 9).Value = 1242.75   # I113: 1191.4 -> 1242.75
$ws.Cells.Item(113, 11).Value = 1242.75   # K113: 1191.4 -> 1242.75
$ws.Cells.Item(113, 13).Value = 927.25   # M113: 978.5999999999999 -> 927.25

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 1909.7391   # H4: 1933.4348 -> 1909.7391
$ws.Cells.Item(4, 9).Value = 1818.6875   # I4: 1909.9333 -> 1818.6875
$ws.Cells.Item(4, 10).Value = 2117.8572   # J4: 1977.5 -> 2117.8572
$ws.Cells.Item(4, 11).Value = 5456.0625   # K4: 5729.7999 -> 5456.0625
$ws.Cells.Item(4, 12).Value = 6353.571599999999   # L4: 5932.5 -> 6353.571599999999
$ws.Cells.Item(4, 13).Value = -5344.0625   # M4: -5617.7999 -> -5344.0625
$ws.Cells.Item(4, 14).Value = -6577.571599999999   # N4: -6156.5 -> -6577.571599999999
$ws.Cells.Item(105, 8).Value = 8600   # H105: 10000 -> 8600
$ws.Cells.Item(105, 10).Value = 8600   # J105: 10000 -> 8600
$ws.Cells.Item(105, 12).Value = 25800   # L105: 30000 -> 25800
$ws.Cells.Item(105, 14).Value = -31042   # N105: -35242 -> -31042
$ws.Cells.Item(117, 8).Value = 730.4   # H117: 739.3333 -> 730.4
$ws.Cells.Item(117, 9).Value = 730.4   # I117: 738.25 -> 730.4
$ws.Cells.Item(117, 10).Value = 0   # J117: 741.5 -> 0
$ws.Cells.Item(117, 11).Value = 2191.2   # K117: 2214.75 -> 2191.2
$ws.Cells.Item(117, 12).Value = 0   # L117: 2224.5 -> 0
$ws.Cells.Item(117, 13).Value = 1250.8   # M117: 1227.25 -> 1250.8
$ws.Cells.Item(117, 14).Value = $null  # N117: clear (was -9108.5)

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(9, 8).Value = 3893.6667   # H9: 403.66666 -> 3893.6667
$ws.Cells.Item(9, 9).Value = 672.4   # I9: 403.66666 -> 672.4
$ws.Cells.Item(9, 10).Value = 20000   # J9: 0 -> 20000
$ws.Cells.Item(9, 11).Value = 672.4   # K9: 403.66666 -> 672.4
$ws.Cells.Item(9, 12).Value = 20000   # L9: 0 -> 20000
$ws.Cells.Item(9, 13).Value = -502.4   # M9: -233.66666 -> -502.4
$ws.Cells.Item(9, 14).Value = -20340   # N9: None -> -20340
$ws.Cells.Item(47, 8).Value = 30000.25   # H47: 30001 -> 30000.25
$ws.Cells.Item(47, 10).Value = 30000.25   # J47: 30001 -> 30000.25
$ws.Cells.Item(47, 12).Value = 30000.25   # L47: 30001 -> 30000.25
$ws.Cells.Item(47, 14).Value = -31136.25   # N47: -31137 -> -31136.25
$ws.Cells.Item(135, 8).Value = 0   # H135: 49999 -> 0
$ws.Cells.Item(135, 10).Value = 0   # J135: 49999 -> 0
$ws.Cells.Item(135, 12).Value = 0   # L135: 49999 -> 0
$ws.Cells.Item(135, 14).Value = $null  # N135: clear (was -60139)

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 3380.6   # H16: 3475.75 -> 3380.6
$ws.Cells.Item(16, 10).Value = 3001   # J16: 3002 -> 3001
$ws.Cells.Item(16, 12).Value = 3001   # L16: 3002 -> 3001
$ws.Cells.Item(16, 14).Value = -3341   # N16: -3342 -> -3341
$ws.Cells.Item(22, 8).Value = 1011.2353   # H22: 1029.4667 -> 1011.2353
$ws.Cells.Item(22, 9).Value = 1054.1   # I22: 1060.2222 -> 1054.1
$ws.Cells.Item(22, 10).Value = 950   # J22: 983.3333 -> 950
$ws.Cells.Item(22, 11).Value = 1054.1   # K22: 1060.2222 -> 1054.1
$ws.Cells.Item(22, 12).Value = 950   # L22: 983.3333 -> 950
$ws.Cells.Item(22, 13).Value = -759.0999999999999   # M22: -765.2221999999999 -> -759.0999999999999
$ws.Cells.Item(22, 14).Value = -1540   # N22: -1573.3333 -> -1540
$ws.Cells.Item(27, 8).Value = 1011.2353   # H27: 1029.4667 -> 1011.2353
$ws.Cells.Item(27, 9).Value = 1054.1   # I27: 1060.2222 -> 1054.1
$ws.Cells.Item(27, 10).Value = 950   # J27: 983.3333 -> 950
$ws.Cells.Item(27, 11).Value = 1054.1   # K27: 1060.2222 -> 1054.1
$ws.Cells.Item(27, 12).Value = 950   # L27: 983.3333 -> 950
$ws.Cells.Item(27, 13).Value = -947.0999999999999   # M27: -953.2221999999999 -> -947.0999999999999
$ws.Cells.Item(61, 8).Value = 404   # H61: 402.66666 -> 404
$ws.Cells.Item(61, 9).Value = 404   # I61: 451.5 -> 404
$ws.Cells.Item(61, 10).Value = 0   # J61: 305 -> 0
$ws.Cells.Item(61, 11).Value = 404   # K61: 451.5 -> 404
$ws.Cells.Item(61, 12).Value = 0   # L61: 305 -> 0
$ws.Cells.Item(61, 13).Value = -202   # M61: -249.5 -> -202
$ws.Cells.Item(61, 14).Value = $null  # N61: clear (was -709)
$ws.Cells.Item(113, 8).Value = 404   # H113: 402.66666 -> 404
$ws.Cells.Item(113, 9).Value = 404   # I113: 451.5 -> 404
$ws.Cells.Item(113, 10).Value = 0   # J113: 305 -> 0
$ws.Cells.Item(113, 11).Value = 404   # K113: 451.5 -> 404
$ws.Cells.Item(113, 12).Value = 0   # L113: 305 -> 0
$ws.Cells.Item(113, 13).Value = 1766   # M113: 1718.5 -> 1766
$ws.Cells.Item(113, 14).Value = $null  # N113: clear (was -4645)

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(113, 8).Value = 834.6667   # H113: 8470 -> 834.6667
$ws.Cells.Item(113, 9).Value = 775.5   # I113: 14963.857 -> 775.5
$ws.Cells.Item(113, 11).Value = 2326.5   # K113: 44891.571 -> 2326.5
$ws.Cells.Item(113, 13).Value = -156.5   # M113: -42721.571 -> -156.5
$ws.Cells.Item(122, 8).Value = 3780.8   # H122: 3817.1667 -> 3780.8
$ws.Cells.Item(122, 9).Value = 3666.3333   # I122: 3749.5 -> 3666.3333
$ws.Cells.Item(122, 11).Value = 10998.9999   # K122: 11248.5 -> 10998.9999
$ws.Cells.Item(122, 13).Value = -8548.999899999999   # M122: -8798.5 -> -8548.999899999999
